$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H2: unchanged text, kept as "Data Update Request"
$ws.Range("H2").Value = "Data Update Request"

# I3: sentences reordered - EPA working note now comes first, scraping note second
$ws.Range("I3").Value = "24/06/2019 - EPA are working on 2014 to 2016 data before sharing. Figures for 2017 will be available later in the year (www.wastereport.ie). Data appears to have been scraped from National Waste Reports which stopped in 2012. Last official update was in Bulletin 2 report for 2013 (published 2014)."

# I4 and I5: replaced with new note about data request for household recycling
$ws.Range("I4").Value = "24/06/2019 - Request for data added to enquiry on data for Waste per Capita (above)"
$ws.Range("I5").Value = "24/06/2019 - Request for data added to enquiry on data for Waste per Capita (above)"

# Update sheet view: selection moves to B6 (and the previous scrolled/frozen
# topLeftCell position is reset to the default view)
$ws.Activate()
$ws.Range("B6").Select()
